# Reorder data rows 2-20 on the active sheet ("Artfynd") according to a
# known row permutation (the underlying cell contents themselves are
# unchanged -- only their row position moves).
#
# old row -> new row
#   2  -> 20
#   3  -> 18
#   4  -> 2
#   5  -> 3
#   6  -> 4
#   7  -> 5
#   8  -> 6
#   9  -> 7
#   10 -> 8
#   11 -> 9
#   12 -> 10
#   13 -> 11
#   14 -> 12
#   15 -> 13
#   16 -> 14
#   17 -> 15
#   18 -> 16
#   19 -> 17
#   20 -> 19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 20
$firstCol = 1
$lastCol = 51   # column AY

# Snapshot every source row (1 x colCount each) before writing anything back,
# since the destination ranges overlap the source ranges.
$srcRows = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
    $srcRows[$r] = $rng.Value2
}

$map = @{
    2  = 20
    3  = 18
    4  = 2
    5  = 3
    6  = 4
    7  = 5
    8  = 6
    9  = 7
    10 = 8
    11 = 9
    12 = 10
    13 = 11
    14 = 12
    15 = 13
    16 = 14
    17 = 15
    18 = 16
    19 = 17
    20 = 19
}

foreach ($oldRow in $map.Keys) {
    $newRow = $map[$oldRow]

    # Columns Y (25) and AA (27) hold date-like text ("2018-06-20").
    # Excel auto-detects that pattern and would silently convert the
    # literal text to a date serial number on assignment, so force those
    # destination cells to text format first to preserve the original
    # string values.
    $ws.Cells.Item($newRow, 25).NumberFormat = "@"
    $ws.Cells.Item($newRow, 27).NumberFormat = "@"

    $rng = $ws.Range($ws.Cells.Item($newRow, $firstCol), $ws.Cells.Item($newRow, $lastCol))
    $rng.Value2 = $srcRows[$oldRow]
}

Write-Output "done"
